{"js": "// Replace arithmetic expressions in each cell of the first table,\n// preserving existing run/paragraph formatting (font, size, alignment).\nconst table = context.document.body.tables.getFirst();\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"31+1=\", newText: \"27+45=\" },\n  { row: 0, col: 1, oldText: \"66+13=\", newText: \"83-24=\" },\n  { row: 0, col: 2, oldText: \"61+7=\", newText: \"58-26=\" },\n  { row: 0, col: 3, oldText: \"5+39=\", newText: \"94-88=\" },\n  { row: 0, col: 4, oldText: \"16+66=\", newText: \"80-24=\" },\n  { row: 1, col: 0, oldText: \"24-17=\", newText: \"98-88=\" },\n  { row: 1, col: 1, oldText: \"96-80=\", newText: \"27+69=\" },\n  { row: 1, col: 2, oldText: \"85-5=\", newText: \"70+3=\" },\n  { row: 1, col: 3, oldText: \"68-1=\", newText: \"39+48=\" },\n  { row: 1, col: 4, oldText: \"22+58=\", newText: \"33-26=\" },\n  { row: 2, col: 0, oldText: \"84-51=\", newText: \"78-12=\" },\n  { row: 2, col: 1, oldText: \"11+77=\", newText: \"62+1=\" },\n  { row: 2, col: 2, oldText: \"43+19=\", newText: \"56-18=\" },\n  { row: 2, col: 3, oldText: \"13+31=\", newText: \"10-9=\" },\n  { row: 2, col: 4, oldText: \"86-63=\", newText: \"43-29=\" },\n  { row: 3, col: 0, oldText: \"56-35=\", newText: \"56-20=\" },\n  { row: 3, col: 1, oldText: \"90-4=\", newText: \"10-7=\" },\n  { row: 3, col: 2, oldText: \"65+11=\", newText: \"91-90=\" },\n  { row: 3, col: 3, oldText: \"95-31=\", newText: \"60+8=\" },\n  { row: 3, col: 4, oldText: \"15-1=\", newText: \"38+53=\" },\n  { row: 4, col: 0, oldText: \"20+11=\", newText: \"47+23=\" },\n  { row: 4, col: 1, oldText: \"21-7=\", newText: \"47+3=\" },\n  { row: 4, col: 2, oldText: \"99-65=\", newText: \"41-22=\" },\n  { row: 4, col: 3, oldText: \"57+7=\", newText: \"43-28=\" },\n  { row: 4, col: 4, oldText: \"84+13=\", newText: \"26-0=\" },\n  { row: 5, col: 0, oldText: \"44-30=\", newText: \"87-20=\" },\n  { row: 5, col: 1, oldText: \"97-28=\", newText: \"49-22=\" },\n  { row: 5, col: 2, oldText: \"62-11=\", newText: \"52-0=\" },\n  { row: 5, col: 3, oldText: \"46+8=\", newText: \"28-0=\" },\n  { row: 5, col: 4, oldText: \"98-70=\", newText: \"67-29=\" },\n  { row: 6, col: 0, oldText: \"36+14=\", newText: \"93-85=\" },\n  { row: 6, col: 1, oldText: \"72+13=\", newText: \"27+26=\" },\n  { row: 6, col: 2, oldText: \"44+6=\", newText: \"56+18=\" },\n  { row: 6, col: 3, oldText: \"74-16=\", newText: \"13+81=\" },\n  { row: 6, col: 4, oldText: \"84-8=\", newText: \"34-21=\" },\n  { row: 7, col: 0, oldText: \"21+33=\", newText: \"67-41=\" },\n  { row: 7, col: 1, oldText: \"70-39=\", newText: \"52-12=\" },\n  { row: 7, col: 2, oldText: \"45+36=\", newText: \"27-20=\" },\n  { row: 7, col: 3, oldText: \"67+32=\", newText: \"76-65=\" },\n  { row: 7, col: 4, oldText: \"77-62=\", newText: \"71-29=\" },\n  { row: 8, col: 0, oldText: \"98-60=\", newText: \"82-59=\" },\n  { row: 8, col: 1, oldText: \"88-51=\", newText: \"3+71=\" },\n  { row: 8, col: 2, oldText: \"20+77=\", newText: \"33+51=\" },\n  { row: 8, col: 3, oldText: \"57+10=\", newText: \"67+11=\" },\n  { row: 8, col: 4, oldText: \"75-32=\", newText: \"51-0=\" },\n  { row: 9, col: 0, oldText: \"12+40=\", newText: \"21+22=\" },\n  { row: 9, col: 1, oldText: \"6-5=\", newText: \"39+49=\" },\n  { row: 9, col: 2, oldText: \"33+21=\", newText: \"85-84=\" },\n  { row: 9, col: 3, oldText: \"40-23=\", newText: \"27+42=\" },\n  { row: 9, col: 4, oldText: \"95-65=\", newText: \"6+53=\" },\n  { row: 10, col: 0, oldText: \"50-48=\", newText: \"59-36=\" },\n  { row: 10, col: 1, oldText: \"2+27=\", newText: \"55-16=\" },\n  { row: 10, col: 2, oldText: \"9+20=\", newText: \"77+8=\" },\n  { row: 10, col: 3, oldText: \"76-70=\", newText: \"45+44=\" },\n  { row: 10, col: 4, oldText: \"14+64=\", newText: \"71-41=\" },\n  { row: 11, col: 0, oldText: \"46+5=\", newText: \"84-64=\" },\n  { row: 11, col: 1, oldText: \"99-24=\", newText: \"89-52=\" },\n  { row: 11, col: 2, oldText: \"83-29=\", newText: \"39+11=\" },\n  { row: 11, col: 3, oldText: \"52-13=\", newText: \"48+6=\" },\n  { row: 11, col: 4, oldText: \"32+59=\", newText: \"29+15=\" },\n  { row: 12, col: 0, oldText: \"11+21=\", newText: \"46-6=\" },\n  { row: 12, col: 1, oldText: \"10+15=\", newText: \"14+5=\" },\n  { row: 12, col: 2, oldText: \"47-33=\", newText: \"14+21=\" },\n  { row: 12, col: 3, oldText: \"79-40=\", newText: \"25+1=\" },\n  { row: 12, col: 4, oldText: \"53+5=\", newText: \"90-53=\" },\n  { row: 13, col: 0, oldText: \"40+21=\", newText: \"60-41=\" },\n  { row: 13, col: 1, oldText: \"8+21=\", newText: \"44+53=\" },\n  { row: 13, col: 2, oldText: \"84-20=\", newText: \"73-67=\" },\n  { row: 13, col: 3, oldText: \"93-6=\", newText: \"90-43=\" },\n  { row: 13, col: 4, oldText: \"95-83=\", newText: \"94+1=\" },\n  { row: 14, col: 0, oldText: \"12+58=\", newText: \"58-10=\" },\n  { row: 14, col: 1, oldText: \"55-17=\", newText: \"21+59=\" },\n  { row: 14, col: 2, oldText: \"15+21=\", newText: \"95-15=\" },\n  { row: 14, col: 3, oldText: \"72-66=\", newText: \"95-36=\" },\n  { row: 14, col: 4, oldText: \"52-36=\", newText: \"8+69=\" },\n  { row: 15, col: 0, oldText: \"53+36=\", newText: \"4+90=\" },\n  { row: 15, col: 1, oldText: \"57-9=\", newText: \"36+21=\" },\n  { row: 15, col: 2, oldText: \"4+26=\", newText: \"74-37=\" },\n  { row: 15, col: 3, oldText: \"56-21=\", newText: \"69-66=\" },\n  { row: 15, col: 4, oldText: \"70-61=\", newText: \"88-2=\" },\n  { row: 16, col: 0, oldText: \"12+50=\", newText: \"95-16=\" },\n  { row: 16, col: 1, oldText: \"21+57=\", newText: \"81-3=\" },\n  { row: 16, col: 2, oldText: \"78-1=\", newText: \"97-70=\" },\n  { row: 16, col: 3, oldText: \"39+38=\", newText: \"59+11=\" },\n  { row: 16, col: 4, oldText: \"7-4=\", newText: \"9+89=\" },\n  { row: 17, col: 0, oldText: \"32+17=\", newText: \"98-53=\" },\n  { row: 17, col: 1, oldText: \"99-21=\", newText: \"72-37=\" },\n  { row: 17, col: 2, oldText: \"94-93=\", newText: \"45+41=\" },\n  { row: 17, col: 3, oldText: \"76+9=\", newText: \"4+80=\" },\n  { row: 17, col: 4, oldText: \"45-23=\", newText: \"25-14=\" },\n  { row: 18, col: 0, oldText: \"62-47=\", newText: \"5+63=\" },\n  { row: 18, col: 1, oldText: \"69-57=\", newText: \"32-6=\" },\n  { row: 18, col: 2, oldText: \"38-10=\", newText: \"91-23=\" },\n  { row: 18, col: 3, oldText: \"30+20=\", newText: \"49-27=\" },\n  { row: 18, col: 4, oldText: \"34+52=\", newText: \"29+49=\" },\n  { row: 19, col: 0, oldText: \"18+71=\", newText: \"95-58=\" },\n  { row: 19, col: 1, oldText: \"48-19=\", newText: \"65-16=\" },\n  { row: 19, col: 2, oldText: \"74+7=\", newText: \"47-41=\" },\n  { row: 19, col: 3, oldText: \"99-8=\", newText: \"99-23=\" },\n  { row: 19, col: 4, oldText: \"40-27=\", newText: \"53-19=\" }\n];\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row, col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  paragraph.load(\"text\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  if (paragraph.text !== oldText) {\n    throw new Error(\n      `Unexpected cell text at row ${row}, col ${col}: expected \"${oldText}\" but found \"${paragraph.text}\"`\n    );\n  }\n\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = '31+1='; NewText = '27+45=' }\n    @{ Row = 1; Col = 2; OldText = '66+13='; NewText = '83-24=' }\n    @{ Row = 1; Col = 3; OldText = '61+7='; NewText = '58-26=' }\n    @{ Row = 1; Col = 4; OldText = '5+39='; NewText = '94-88=' }\n    @{ Row = 1; Col = 5; OldText = '16+66='; NewText = '80-24=' }\n    @{ Row = 2; Col = 1; OldText = '24-17='; NewText = '98-88=' }\n    @{ Row = 2; Col = 2; OldText = '96-80='; NewText = '27+69=' }\n    @{ Row = 2; Col = 3; OldText = '85-5='; NewText = '70+3=' }\n    @{ Row = 2; Col = 4; OldText = '68-1='; NewText = '39+48=' }\n    @{ Row = 2; Col = 5; OldText = '22+58='; NewText = '33-26=' }\n    @{ Row = 3; Col = 1; OldText = '84-51='; NewText = '78-12=' }\n    @{ Row = 3; Col = 2; OldText = '11+77='; NewText = '62+1=' }\n    @{ Row = 3; Col = 3; OldText = '43+19='; NewText = '56-18=' }\n    @{ Row = 3; Col = 4; OldText = '13+31='; NewText = '10-9=' }\n    @{ Row = 3; Col = 5; OldText = '86-63='; NewText = '43-29=' }\n    @{ Row = 4; Col = 1; OldText = '56-35='; NewText = '56-20=' }\n    @{ Row = 4; Col = 2; OldText = '90-4='; NewText = '10-7=' }\n    @{ Row = 4; Col = 3; OldText = '65+11='; NewText = '91-90=' }\n    @{ Row = 4; Col = 4; OldText = '95-31='; NewText = '60+8=' }\n    @{ Row = 4; Col = 5; OldText = '15-1='; NewText = '38+53=' }\n    @{ Row = 5; Col = 1; OldText = '20+11='; NewText = '47+23=' }\n    @{ Row = 5; Col = 2; OldText = '21-7='; NewText = '47+3=' }\n    @{ Row = 5; Col = 3; OldText = '99-65='; NewText = '41-22=' }\n    @{ Row = 5; Col = 4; OldText = '57+7='; NewText = '43-28=' }\n    @{ Row = 5; Col = 5; OldText = '84+13='; NewText = '26-0=' }\n    @{ Row = 6; Col = 1; OldText = '44-30='; NewText = '87-20=' }\n    @{ Row = 6; Col = 2; OldText = '97-28='; NewText = '49-22=' }\n    @{ Row = 6; Col = 3; OldText = '62-11='; NewText = '52-0=' }\n    @{ Row = 6; Col = 4; OldText = '46+8='; NewText = '28-0=' }\n    @{ Row = 6; Col = 5; OldText = '98-70='; NewText = '67-29=' }\n    @{ Row = 7; Col = 1; OldText = '36+14='; NewText = '93-85=' }\n    @{ Row = 7; Col = 2; OldText = '72+13='; NewText = '27+26=' }\n    @{ Row = 7; Col = 3; OldText = '44+6='; NewText = '56+18=' }\n    @{ Row = 7; Col = 4; OldText = '74-16='; NewText = '13+81=' }\n    @{ Row = 7; Col = 5; OldText = '84-8='; NewText = '34-21=' }\n    @{ Row = 8; Col = 1; OldText = '21+33='; NewText = '67-41=' }\n    @{ Row = 8; Col = 2; OldText = '70-39='; NewText = '52-12=' }\n    @{ Row = 8; Col = 3; OldText = '45+36='; NewText = '27-20=' }\n    @{ Row = 8; Col = 4; OldText = '67+32='; NewText = '76-65=' }\n    @{ Row = 8; Col = 5; OldText = '77-62='; NewText = '71-29=' }\n    @{ Row = 9; Col = 1; OldText = '98-60='; NewText = '82-59=' }\n    @{ Row = 9; Col = 2; OldText = '88-51='; NewText = '3+71=' }\n    @{ Row = 9; Col = 3; OldText = '20+77='; NewText = '33+51=' }\n    @{ Row = 9; Col = 4; OldText = '57+10='; NewText = '67+11=' }\n    @{ Row = 9; Col = 5; OldText = '75-32='; NewText = '51-0=' }\n    @{ Row = 10; Col = 1; OldText = '12+40='; NewText = '21+22=' }\n    @{ Row = 10; Col = 2; OldText = '6-5='; NewText = '39+49=' }\n    @{ Row = 10; Col = 3; OldText = '33+21='; NewText = '85-84=' }\n    @{ Row = 10; Col = 4; OldText = '40-23='; NewText = '27+42=' }\n    @{ Row = 10; Col = 5; OldText = '95-65='; NewText = '6+53=' }\n    @{ Row = 11; Col = 1; OldText = '50-48='; NewText = '59-36=' }\n    @{ Row = 11; Col = 2; OldText = '2+27='; NewText = '55-16=' }\n    @{ Row = 11; Col = 3; OldText = '9+20='; NewText = '77+8=' }\n    @{ Row = 11; Col = 4; OldText = '76-70='; NewText = '45+44=' }\n    @{ Row = 11; Col = 5; OldText = '14+64='; NewText = '71-41=' }\n    @{ Row = 12; Col = 1; OldText = '46+5='; NewText = '84-64=' }\n    @{ Row = 12; Col = 2; OldText = '99-24='; NewText = '89-52=' }\n    @{ Row = 12; Col = 3; OldText = '83-29='; NewText = '39+11=' }\n    @{ Row = 12; Col = 4; OldText = '52-13='; NewText = '48+6=' }\n    @{ Row = 12; Col = 5; OldText = '32+59='; NewText = '29+15=' }\n    @{ Row = 13; Col = 1; OldText = '11+21='; NewText = '46-6=' }\n    @{ Row = 13; Col = 2; OldText = '10+15='; NewText = '14+5=' }\n    @{ Row = 13; Col = 3; OldText = '47-33='; NewText = '14+21=' }\n    @{ Row = 13; Col = 4; OldText = '79-40='; NewText = '25+1=' }\n    @{ Row = 13; Col = 5; OldText = '53+5='; NewText = '90-53=' }\n    @{ Row = 14; Col = 1; OldText = '40+21='; NewText = '60-41=' }\n    @{ Row = 14; Col = 2; OldText = '8+21='; NewText = '44+53=' }\n    @{ Row = 14; Col = 3; OldText = '84-20='; NewText = '73-67=' }\n    @{ Row = 14; Col = 4; OldText = '93-6='; NewText = '90-43=' }\n    @{ Row = 14; Col = 5; OldText = '95-83='; NewText = '94+1=' }\n    @{ Row = 15; Col = 1; OldText = '12+58='; NewText = '58-10=' }\n    @{ Row = 15; Col = 2; OldText = '55-17='; NewText = '21+59=' }\n    @{ Row = 15; Col = 3; OldText = '15+21='; NewText = '95-15=' }\n    @{ Row = 15; Col = 4; OldText = '72-66='; NewText = '95-36=' }\n    @{ Row = 15; Col = 5; OldText = '52-36='; NewText = '8+69=' }\n    @{ Row = 16; Col = 1; OldText = '53+36='; NewText = '4+90=' }\n    @{ Row = 16; Col = 2; OldText = '57-9='; NewText = '36+21=' }\n    @{ Row = 16; Col = 3; OldText = '4+26='; NewText = '74-37=' }\n    @{ Row = 16; Col = 4; OldText = '56-21='; NewText = '69-66=' }\n    @{ Row = 16; Col = 5; OldText = '70-61='; NewText = '88-2=' }\n    @{ Row = 17; Col = 1; OldText = '12+50='; NewText = '95-16=' }\n    @{ Row = 17; Col = 2; OldText = '21+57='; NewText = '81-3=' }\n    @{ Row = 17; Col = 3; OldText = '78-1='; NewText = '97-70=' }\n    @{ Row = 17; Col = 4; OldText = '39+38='; NewText = '59+11=' }\n    @{ Row = 17; Col = 5; OldText = '7-4='; NewText = '9+89=' }\n    @{ Row = 18; Col = 1; OldText = '32+17='; NewText = '98-53=' }\n    @{ Row = 18; Col = 2; OldText = '99-21='; NewText = '72-37=' }\n    @{ Row = 18; Col = 3; OldText = '94-93='; NewText = '45+41=' }\n    @{ Row = 18; Col = 4; OldText = '76+9='; NewText = '4+80=' }\n    @{ Row = 18; Col = 5; OldText = '45-23='; NewText = '25-14=' }\n    @{ Row = 19; Col = 1; OldText = '62-47='; NewText = '5+63=' }\n    @{ Row = 19; Col = 2; OldText = '69-57='; NewText = '32-6=' }\n    @{ Row = 19; Col = 3; OldText = '38-10='; NewText = '91-23=' }\n    @{ Row = 19; Col = 4; OldText = '30+20='; NewText = '49-27=' }\n    @{ Row = 19; Col = 5; OldText = '34+52='; NewText = '29+49=' }\n    @{ Row = 20; Col = 1; OldText = '18+71='; NewText = '95-58=' }\n    @{ Row = 20; Col = 2; OldText = '48-19='; NewText = '65-16=' }\n    @{ Row = 20; Col = 3; OldText = '74+7='; NewText = '47-41=' }\n    @{ Row = 20; Col = 4; OldText = '99-8='; NewText = '99-23=' }\n    @{ Row = 20; Col = 5; OldText = '40-27='; NewText = '53-19=' }\n)\n\nforeach ($item in $replacements) {\n    $cell = $table.Cell($item.Row, $item.Col)\n    $r = $cell.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    if ($r.Text -ne $item.OldText) {\n        throw \"Unexpected cell text at row $($item.Row), col $($item.Col): expected $($item.OldText) but found $($r.Text)\"\n    }\n    $r.Text = $item.NewText\n}"}
